$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) for rows 2-5 from 2023-10-13 (45212) to 2023-10-22 (45221)
$ws.Range("C2").Value = 45221
$ws.Range("C3").Value = 45221
$ws.Range("C4").Value = 45221
$ws.Range("C5").Value = 45221
